$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (column A date style, borders, etc.) from the last
# existing row down through the new rows in one shot.
$ws.Range("A464:D464").Copy()
$ws.Range("A465:D491").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new daily figures (9 dic 2021 -> 5 gen 2022, "fino a 6
# gennaio 2022" per il commit) riportati nel foglio.
$ws.Cells.Item(465, 1).Value = 44539
$ws.Cells.Item(465, 2).Value = 7
$ws.Cells.Item(465, 3).Value = 43
$ws.Cells.Item(465, 4).Value = 125.1236687423616

$ws.Cells.Item(466, 1).Value = 44540
$ws.Cells.Item(466, 2).Value = 11
$ws.Cells.Item(466, 3).Value = 52
$ws.Cells.Item(466, 4).Value = 151.3123435954141

$ws.Cells.Item(467, 1).Value = 44541
$ws.Cells.Item(467, 2).Value = 2
$ws.Cells.Item(467, 3).Value = 41
$ws.Cells.Item(467, 4).Value = 119.3039632194611

$ws.Cells.Item(468, 1).Value = 44542
$ws.Cells.Item(468, 2).Value = 24
$ws.Cells.Item(468, 3).Value = 63
$ws.Cells.Item(468, 4).Value = 183.320723971367

$ws.Cells.Item(469, 1).Value = 44543
$ws.Cells.Item(469, 2).Value = 12
$ws.Cells.Item(469, 3).Value = 71
$ws.Cells.Item(469, 4).Value = 206.5995460629692

$ws.Cells.Item(470, 1).Value = 44544
$ws.Cells.Item(470, 2).Value = 14
$ws.Cells.Item(470, 3).Value = 76
$ws.Cells.Item(470, 4).Value = 221.1488098702206

$ws.Cells.Item(471, 1).Value = 44545
$ws.Cells.Item(471, 2).Value = 6
$ws.Cells.Item(471, 3).Value = 76
$ws.Cells.Item(471, 4).Value = 221.1488098702206

$ws.Cells.Item(472, 1).Value = 44546
$ws.Cells.Item(472, 2).Value = 14
$ws.Cells.Item(472, 3).Value = 83
$ws.Cells.Item(472, 4).Value = 241.5177792003724

$ws.Cells.Item(473, 1).Value = 44547
$ws.Cells.Item(473, 2).Value = 6
$ws.Cells.Item(473, 3).Value = 78
$ws.Cells.Item(473, 4).Value = 226.9685153931211

$ws.Cells.Item(474, 1).Value = 44548
$ws.Cells.Item(474, 2).Value = 14
$ws.Cells.Item(474, 3).Value = 90
$ws.Cells.Item(474, 4).Value = 261.8867485305244

$ws.Cells.Item(475, 1).Value = 44550
$ws.Cells.Item(475, 2).Value = 16
$ws.Cells.Item(475, 3).Value = 82
$ws.Cells.Item(475, 4).Value = 238.6079264389222

$ws.Cells.Item(476, 1).Value = 44551
$ws.Cells.Item(476, 2).Value = 24
$ws.Cells.Item(476, 3).Value = 94
$ws.Cells.Item(476, 4).Value = 273.5261595763254

$ws.Cells.Item(477, 1).Value = 44552
$ws.Cells.Item(477, 2).Value = 10
$ws.Cells.Item(477, 3).Value = 90
$ws.Cells.Item(477, 4).Value = 261.8867485305244

$ws.Cells.Item(478, 1).Value = 44553
$ws.Cells.Item(478, 2).Value = 28
$ws.Cells.Item(478, 3).Value = 112
$ws.Cells.Item(478, 4).Value = 325.9035092824303

$ws.Cells.Item(479, 1).Value = 44554
$ws.Cells.Item(479, 2).Value = 16
$ws.Cells.Item(479, 3).Value = 114
$ws.Cells.Item(479, 4).Value = 331.7232148053309

$ws.Cells.Item(480, 1).Value = 44555
$ws.Cells.Item(480, 2).Value = 57
$ws.Cells.Item(480, 3).Value = 165
$ws.Cells.Item(480, 4).Value = 480.1257056392947

$ws.Cells.Item(481, 1).Value = 44556
$ws.Cells.Item(481, 2).Value = 49
$ws.Cells.Item(481, 3).Value = 200
$ws.Cells.Item(481, 4).Value = 581.9705522900541

$ws.Cells.Item(482, 1).Value = 44557
$ws.Cells.Item(482, 2).Value = 48
$ws.Cells.Item(482, 3).Value = 232
$ws.Cells.Item(482, 4).Value = 675.0858406564628

$ws.Cells.Item(483, 1).Value = 44558
$ws.Cells.Item(483, 2).Value = 15
$ws.Cells.Item(483, 3).Value = 223
$ws.Cells.Item(483, 4).Value = 648.8971658034103

$ws.Cells.Item(484, 1).Value = 44559
$ws.Cells.Item(484, 2).Value = 63
$ws.Cells.Item(484, 3).Value = 276
$ws.Cells.Item(484, 4).Value = 803.1193621602747

$ws.Cells.Item(485, 1).Value = 44560
$ws.Cells.Item(485, 2).Value = 110
$ws.Cells.Item(485, 3).Value = 358
$ws.Cells.Item(485, 4).Value = 1041.727288599197

$ws.Cells.Item(486, 1).Value = 44561
$ws.Cells.Item(486, 2).Value = 141
$ws.Cells.Item(486, 3).Value = 483
$ws.Cells.Item(486, 4).Value = 1405.458883780481

$ws.Cells.Item(487, 1).Value = 44562
$ws.Cells.Item(487, 2).Value = 92
$ws.Cells.Item(487, 3).Value = 518
$ws.Cells.Item(487, 4).Value = 1507.30373043124

$ws.Cells.Item(488, 1).Value = 44563
$ws.Cells.Item(488, 2).Value = 38
$ws.Cells.Item(488, 3).Value = 507
$ws.Cells.Item(488, 4).Value = 1475.295350055287

$ws.Cells.Item(489, 1).Value = 44564
$ws.Cells.Item(489, 2).Value = 34
$ws.Cells.Item(489, 3).Value = 493
$ws.Cells.Item(489, 4).Value = 1434.557411394983

$ws.Cells.Item(490, 1).Value = 44565
$ws.Cells.Item(490, 2).Value = 12
$ws.Cells.Item(490, 3).Value = 490
$ws.Cells.Item(490, 4).Value = 1425.827853110633

$ws.Cells.Item(491, 1).Value = 44566
$ws.Cells.Item(491, 2).Value = 37
$ws.Cells.Item(491, 3).Value = 464
$ws.Cells.Item(491, 4).Value = 1350.171681312926
